$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeImage($range, $oldName, $newName) {
    if ($range.InlineShapes.Count -gt 0) {
        for ($j = 1; $j -le $range.InlineShapes.Count; $j++) {
            $inline = $range.InlineShapes.Item($j)
            $floatShape = $inline.ConvertToShape()
            if ($floatShape.Name -eq $oldName) {
                $floatShape.Name = $newName
            }
            $floatShape.ConvertToInlineShape() | Out-Null
        }
    }
}

# Footers: the Pearson logo picture's docPr/cNvPr "name" goes from
# image2.png -> image1.png (both footer slots use the same picture name).
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        Rename-InlineShapeImage $ftr.Range "image2.png" "image1.png"
    }
}

# Headers: the BTec logo picture's docPr/cNvPr "name" goes from
# image1.jpg -> image2.jpg.
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        Rename-InlineShapeImage $hdr.Range "image1.jpg" "image2.jpg"
    }
}
